$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13:50 down to 14:51
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with data
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(13, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
$ws.Cells.Item(13, 4).Value = 44498
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = 100114001
$ws.Cells.Item(13, 7).Value = "Papa"
$ws.Cells.Item(13, 8).Value = "Rosara"
$ws.Cells.Item(13, 9).Value = "1a nueva(o)"
$ws.Cells.Item(13, 10).Value = 1000
$ws.Cells.Item(13, 11).Value = 14000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 14500
$ws.Cells.Item(13, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Región del Maule"
$ws.Cells.Item(13, 16).Value = 580
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
